$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.428.44'
$ws.Range('E2').Value = '  -4.80%  '
$ws.Range('D3').Value = '3.114.73'
$ws.Range('E3').Value = '  -5.82%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''520.43'
$ws.Range('E5').Value = '  -6.86%  '
$ws.Range('D6').Value = '''134.08'
$ws.Range('E6').Value = '  -5.49%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '3.112.82'
$ws.Range('E8').Value = '  -5.92%  '
$ws.Range('E9').Value = '  -6.72%  '
$ws.Range('E10').Value = '  -8.72%  '
$ws.Range('E11').Value = '  -7.96%  '
$ws.Range('E12').Value = '  -5.90%  '
$ws.Range('D13').Value = '3.651.90'
$ws.Range('E13').Value = '  -5.84%  '
$ws.Range('E14').Value = '  -2.52%  '
$ws.Range('D15').Value = '''25.33'
$ws.Range('E15').Value = '  -5.49%  '
$ws.Range('D16').Value = '3.110.40'
$ws.Range('E16').Value = '  -5.89%  '
$ws.Range('D17').Value = '57.397.54'
$ws.Range('E17').Value = '  -4.87%  '
$ws.Range('D18').Value = '''0.0000151'
$ws.Range('E18').Value = '  -8.78%  '
$ws.Range('D19').Value = '''5.77'
$ws.Range('E19').Value = '  -6.53%  '
$ws.Range('D20').Value = '''12.95'
$ws.Range('E20').Value = '  -10.39%  '
$ws.Range('D21').Value = '''7.95'
$ws.Range('E21').Value = '  -8.02%  '
$ws.Range('D22').Value = '''340.92'
$ws.Range('E22').Value = '  -9.11%  '
$ws.Range('D23').Value = '''1.00'
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').Value = '''68.10'
$ws.Range('E24').Value = '  -8.23%  '
$ws.Range('D25').Value = '''0.502'
$ws.Range('E25').Value = '  -7.44%  '
$ws.Range('D26').Value = '3.244.43'
$ws.Range('E26').Value = '  -5.84%  '
$ws.Range('E27').Value = '  -4.17%  '
$ws.Range('D28').Value = '''0.997'
$ws.Range('E28').Value = '  -0.37%  '
$ws.Range('D29').Value = '0.0₃0936'
$ws.Range('E29').Value = '  -9.19%  '
$ws.Range('E30').Value = '  -0.13%  '
$ws.Range('D31').Value = '''6.73'
$ws.Range('E31').Value = '  -7.12%  '
$ws.Range('E32').Value = '  -8.59%  '
$ws.Range('D33').Value = '''6.87'
$ws.Range('E33').Value = '  -10.03%  '
$ws.Range('D34').Value = '''21.42'
$ws.Range('E34').Value = '  -5.00%  '
$ws.Range('E35').Value = '  -3.75%  '
$ws.Range('D36').Value = '''157.81'
$ws.Range('E36').Value = '  -4.94%  '
$ws.Range('D37').Value = '''4.74'
$ws.Range('E37').Value = '  -7.67%  '
$ws.Range('D38').Value = '''6.14'
$ws.Range('E38').Value = '  -8.63%  '
$ws.Range('D39').Value = '''1.37'
$ws.Range('E39').Value = '  -10.85%  '
$ws.Range('D40').Value = '''25.12'
$ws.Range('E40').Value = '  -6.22%  '
$ws.Range('E41').Value = '  -6.93%  '
$ws.Range('D42').Value = '3.144.16'
$ws.Range('E42').Value = '  -5.85%  '
$ws.Range('D43').Value = '''40.25'
$ws.Range('E43').Value = '  -4.13%  '
$ws.Range('D44').Value = '''0.680'
$ws.Range('E44').Value = '  -9.53%  '
$ws.Range('E45').Value = '  -5.02%  '
$ws.Range('E46').Value = '  -7.33%  '
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('E48').Value = '  -9.61%  '
$ws.Range('D49').Value = '2.246.64'
$ws.Range('E49').Value = '  -5.35%  '
$ws.Range('E50').Value = '  -5.61%  '
$ws.Range('D51').Value = '''19.94'
$ws.Range('E51').Value = '  -6.34%  '
